# feat: add 2022-Q1 data
#
# The original "总计" (Total) sheet is repurposed into the new "2022-Q1"
# sheet (keeping its sheetId / relationship slot, matching how the source
# workbook was actually authored), and a brand-new "总计" sheet is appended
# at the end of the workbook with the combined/refreshed totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: turn the existing "总计" sheet into the new "2022-Q1" sheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Drop the old second data row (used to hold the 2021-Q2 total); the
# 2022-Q1 sheet only needs a single fund row.
$q1.Rows.Item(3).Delete()

# Stretch the header formatting (bold/border/center) from D1 into the
# newly needed E1:H1 header cells.
$q1.Range("D1").Copy($q1.Range("E1:H1"))

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Fund data row (values kept as text, matching the workbook's convention
# for these numeric-looking figures).
$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "003366"
$q1.Range("B2").Style = "Normal"

$q1.Range("C2").Value = "浙商汇金中证转型成长指数"

$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "0.09"
$q1.Range("D2").Style = "Normal"

$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "93.88"
$q1.Range("E2").Style = "Normal"

$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "1.27"
$q1.Range("F2").Style = "Normal"

$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.0011"
$q1.Range("G2").Style = "Normal"

$q1.Range("H2").Value = 2

# ---------------------------------------------------------------------
# Step 2: append a brand-new "总计" sheet with the refreshed totals
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add($null, $lastSheet)
$total.Name = "总计"

# Match the page setup used by the other sheets.
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36
$total.Outline.SummaryRow = 1
$total.Outline.SummaryColumn = 1

# Borrow the bold/bordered header style (and the matching "index" column
# style) from the 2021-Q3 sheet.
$styleSrc = $wb.Worksheets.Item("2021-Q3")
$styleSrc.Range("B1:D1").Copy($total.Range("B1:D1"))
$styleSrc.Range("A2").Copy($total.Range("A2"))
$styleSrc.Range("A2").Copy($total.Range("A3"))
$styleSrc.Range("A2").Copy($total.Range("A4"))

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q3"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.01

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q2"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.01

# ---------------------------------------------------------------------
# Step 3: restore original active sheet/tab
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
